# This script reshapes the "Data_Full" worksheet:
#   - Inserts a new column A holding a 0-based numeric "segments" index
#   - Shifts the old columns A-F one column to the right (B-G)
#   - Adds a "segments" header label in the new B1 header cell
#   - The old segment-name column (now column B) loses the bold/border
#     header-like style that it used to carry, while the new numeric
#     index column (A) picks that style up instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 20
$firstDataRow = 2

# --- 1. Read all existing values before we start overwriting anything ---
$oldA = @{}   # segment name strings (A2:A20)
$oldB = @{}
$oldC = @{}
$oldD = @{}
$oldE = @{}
$oldF = @{}

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $oldA[$r] = $ws.Range("A$r").Value2
    $oldB[$r] = $ws.Range("B$r").Value2
    $oldC[$r] = $ws.Range("C$r").Value2
    $oldD[$r] = $ws.Range("D$r").Value2
    $oldE[$r] = $ws.Range("E$r").Value2
    $oldF[$r] = $ws.Range("F$r").Value2
}

$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2
$oldF1 = $ws.Range("F1").Value2

# --- 2. Prepare the style that the new index column (A) should use.
#        Column B1 (a header cell) already carries the bold/border style
#        ("style 1") that we want column A's data cells to use too. The
#        new header cell G1 is brand new territory (outside the old
#        A1:F20 used range) and therefore starts out unstyled, so make
#        sure it also picks up that same header style explicitly. ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G1").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Shift header row one column to the right and add the new label ---
$ws.Range("G1").Value = $oldF1
$ws.Range("F1").Value = $oldE1
$ws.Range("E1").Value = $oldD1
$ws.Range("D1").Value = $oldC1
$ws.Range("C1").Value = $oldB1
$ws.Range("B1").Value = "segments"

# --- 4. Shift each data row one column to the right and fill in the
#        new 0-based numeric index in column A. The old segment-name
#        text (previously in column A) moves into column B, which
#        keeps no special style (matching the rest of the data cells). ---
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $ws.Range("G$r").Value = $oldF[$r]
    $ws.Range("F$r").Value = $oldE[$r]
    $ws.Range("E$r").Value = $oldD[$r]
    $ws.Range("D$r").Value = $oldC[$r]
    $ws.Range("C$r").Value = $oldB[$r]
    $ws.Range("B$r").Value = $oldA[$r]
    $ws.Range("A$r").Value = $r - $firstDataRow
}

Write-Output "Reshape complete"
